$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("2010 and 2010-18")
$ws2 = $wb.Worksheets.Item("2000-09")

# --- New data rows on sheet "2010 and 2010-18" (rows 55 and 56) ---
# Row 55: Baseline 2010-18_C167 with m_n = 0.15
$ws1.Cells.Item(55, 1).Value = "CW3M"
$ws1.Cells.Item(55, 2).Value = "Baseline 2010-18_C167 with m_n = 0.15"
$ws1.Cells.Item(55, 3).Value = "2010-18"
$ws1.Cells.Item(55, 4).Value = 1207.0222438888889
$ws1.Cells.Item(55, 5).Value = 1901.5157334444443
$ws1.Cells.Item(55, 6).Value = 0.97970299999999988
$ws1.Cells.Item(55, 7).Value = 280.33542888888883
$ws1.Cells.Item(55, 8).Value = 9.775355222222224
$ws1.Cells.Item(55, 9).Value = 5.3316344444444441
$ws1.Cells.Item(55, 10).Value = 8.145128999999999
$ws1.Cells.Item(55, 11).Value = 645.95907266666666
$ws1.Cells.Item(55, 12).Value = 83.47062044444445
$ws1.Cells.Item(55, 13).Value = 1460.4185112222222
$ws1.Cells.Item(55, 14).Value = 1206.2352837777776
$ws1.Cells.Item(55, 15).Value = 4662.5755209999998
$ws1.Cells.Item(55, 16).Value = 27227.338324888889
$ws1.Cells.Item(55, 17).Value = -0.73148200000000008
$ws1.Cells.Item(55, 18).Value = -0.00020655555555555556

# Row 56: Baseline 2010-18_C167 with m_n = 0.04
$ws1.Cells.Item(56, 1).Value = "CW3M"
$ws1.Cells.Item(56, 2).Value = "Baseline 2010-18_C167 with m_n = 0.04"
$ws1.Cells.Item(56, 3).Value = "2010-18"
$ws1.Cells.Item(56, 4).Value = 1206.5233695555557
$ws1.Cells.Item(56, 5).Value = 1901.5157334444443
$ws1.Cells.Item(56, 6).Value = 0.97970299999999988
$ws1.Cells.Item(56, 7).Value = 280.33542888888883
$ws1.Cells.Item(56, 8).Value = 9.775355222222224
$ws1.Cells.Item(56, 9).Value = 5.3531247777777775
$ws1.Cells.Item(56, 10).Value = 8.145128999999999
$ws1.Cells.Item(56, 11).Value = 645.94818811111122
$ws1.Cells.Item(56, 12).Value = 83.47062044444445
$ws1.Cells.Item(56, 13).Value = 1460.5092637777777
$ws1.Cells.Item(56, 14).Value = 1205.5020886666666
$ws1.Cells.Item(56, 15).Value = 4662.5708008888896
$ws1.Cells.Item(56, 16).Value = 27227.338324888889
$ws1.Cells.Item(56, 17).Value = -0.90742522222222222
$ws1.Cells.Item(56, 18).Value = -0.00025755555555555558

# --- Number formats matching the rest of the table ---
$ws1.Range("D55:N56").NumberFormat = "0.00"
$ws1.Range("O55:P56").NumberFormat = "0"
$ws1.Range("Q55:Q56").NumberFormat = "0.00"
$ws1.Range("R55:R56").NumberFormat = "0.000000"

# --- Column B wraps text like the other scenario-name cells ---
$ws1.Range("B55:B56").WrapText = $true
$ws1.Rows.Item(55).RowHeight = 28.8
$ws1.Rows.Item(56).RowHeight = 28.8

# --- Sheet view / window state changes ---
# The first sheet becomes the tab that is selected/active again, and its
# frozen top row is re-anchored (freeze still covers just row 1) before the
# final selection lands on the new scenario's name cell.
$ws1.Select()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("B57").Select()

$ws2.Range("E3").Select()
$ws1.Select()
